$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "99.068.43"
$ws.Range("E2").Value = "  +2.39%  "
$ws.Range("D3").Value = "3.410.74"
$ws.Range("E3").Value = "  +8.88%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "261.45"
$ws.Range("E5").Value = "  +8.67%  "
$ws.Range("D6").Value = "636.21"
$ws.Range("E6").Value = "  +4.36%  "
$ws.Range("E7").Value = "  +26.50%  "
$ws.Range("D8").Value = "0.397"
$ws.Range("E8").Value = "  +2.17%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "0.892"
$ws.Range("E10").Value = "  +11.60%  "
$ws.Range("D11").Value = "3.410.23"
$ws.Range("E11").Value = "  +9.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.200"
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("D13").Value = "98.593.48"
$ws.Range("E13").Value = "  +2.81%  "
$ws.Range("D14").Value = "36.59"
$ws.Range("E14").Value = "  +6.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000252"
$ws.Range("E15").Value = "  +4.18%  "
$ws.Range("D16").Value = "4.044.00"
$ws.Range("E16").Value = "  +9.29%  "
$ws.Range("E17").Value = "  +4.04%  "
$ws.Range("D18").Value = "3.389.96"
$ws.Range("E18").Value = "  +9.16%  "
$ws.Range("D19").Value = "3.63"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("D20").Value = "15.33"
$ws.Range("E20").Value = "  +5.45%  "
$ws.Range("D21").Value = "495.98"
$ws.Range("E21").Value = "  +2.91%  "
$ws.Range("D22").Value = "6.24"
$ws.Range("E22").Value = "  +8.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000215"
$ws.Range("E23").Value = "  +9.46%  "
$ws.Range("D24").Value = "9.49"
$ws.Range("E24").Value = "  +7.49%  "
$ws.Range("D25").Value = "5.78"
$ws.Range("E25").Value = "  +3.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "89.40"
$ws.Range("E26").Value = "  +4.56%  "
$ws.Range("E27").Value = "  +3.33%  "
$ws.Range("D28").Value = "3.519.11"
$ws.Range("E28").Value = "  +7.80%  "
$ws.Range("D29").Value = "0.288"
$ws.Range("E29").Value = "  +20.41%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "0.194"
$ws.Range("E31").Value = "  +9.65%  "
$ws.Range("D32").Value = "0.134"
$ws.Range("E32").Value = "  +6.10%  "
$ws.Range("D33").Value = "9.76"
$ws.Range("E33").Value = "  +6.97%  "
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("D35").Value = "28.18"
$ws.Range("E35").Value = "  +7.00%  "
$ws.Range("D36").Value = "7.45"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").Value = "0.151"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.00"
$ws.Range("E38").Value = "  +6.48%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "0.478"
$ws.Range("E39").Value = "  +8.27%  "
$ws.Range("B40").Value = "Bittensor"
$ws.Range("C40").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D40").Value = "510.04"
$ws.Range("E40").Value = "  +2.79%  "
$ws.Range("D41").Value = "24.88"
$ws.Range("E41").Value = "  +2.91%  "
$ws.Range("D42").Value = "3.84"
$ws.Range("E42").Value = "  +5.45%  "
$ws.Range("E43").Value = "  +3.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.40"
$ws.Range("E44").Value = "  +5.37%  "
$ws.Range("D45").Value = "0.794"
$ws.Range("E45").Value = "  +13.15%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "160.31"
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").Value = "1.96"
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("D49").Value = "4.73"
$ws.Range("E49").Value = "  +8.32%  "
$ws.Range("D50").Value = "46.82"
$ws.Range("E50").Value = "  +6.31%  "
$ws.Range("D51").Value = "0.825"
$ws.Range("E51").Value = "  +12.06%  "
